# Auto-generated Excel COM-interop script
# Re-applies the periodic cryptos.xlsx price / 1h-volume refresh
# (GitHub Actions bot commit) captured by the supplied OOXML diff.
#
# Column layout: A=index, B=Coin, C=Link, D=Price, E=Volume(1h)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.940.83'
$ws.Range('E2').Value = '  +1.15%  '
$ws.Range('D3').Value = '2.452.00'
$ws.Range('E3').Value = '  +2.36%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''563.49'
$ws.Range('E5').Value = '  +1.00%  '
$ws.Range('D6').Value = '''141.98'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  -0.71%  '
$ws.Range('D9').Value = '2.452.03'
$ws.Range('E9').Value = '  +2.49%  '
$ws.Range('E10').Value = '  -0.24%  '
$ws.Range('E11').Value = '  -0.24%  '
$ws.Range('E12').Value = '  +1.68%  '
$ws.Range('D13').Value = '''0.353'
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').Value = '''26.93'
$ws.Range('E14').Value = '  +4.29%  '
$ws.Range('D15').Value = '2.891.44'
$ws.Range('E15').Value = '  +2.45%  '
$ws.Range('D16').Value = '62.795.81'
$ws.Range('E16').Value = '  +1.06%  '
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').Value = '2.452.19'
$ws.Range('E18').Value = '  +2.93%  '
$ws.Range('D19').Value = '''11.19'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').Value = '''338.62'
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').Value = '''4.26'
$ws.Range('E21').Value = '  +0.78%  '
$ws.Range('D22').Value = '''6.73'
$ws.Range('E22').Value = '  -1.60%  '
$ws.Range('D23').Value = '''0.999'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = '''65.35'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').Value = '''0.170'
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('E28').Value = '  +4.16%  '
$ws.Range('D29').Value = '''8.03'
$ws.Range('E29').Value = '  -3.54%  '
$ws.Range('D30').Value = '''6.73'
$ws.Range('E30').Value = '  +6.26%  '
$ws.Range('E31').Value = '  +1.82%  '
$ws.Range('E32').Value = '  +2.57%  '
$ws.Range('D33').Value = '''176.48'
$ws.Range('E33').Value = '  +3.21%  '
$ws.Range('E34').Value = '  +7.41%  '
$ws.Range('D35').Value = '''384.51'
$ws.Range('E35').Value = '  +9.94%  '
$ws.Range('D36').Value = '''0.395'
$ws.Range('E36').Value = '  +0.84%  '
$ws.Range('D37').Value = '''18.69'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('E38').Value = '  +0.02%  '
$ws.Range('E39').Value = '  -4.54%  '
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').Value = '''1.72'
$ws.Range('E40').Value = '  +6.75%  '
$ws.Range('B41').Value = 'FirstDigitalUSD'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''0.999'
$ws.Range('E41').Value = '  -0.19%  '
$ws.Range('D42').Value = '''39.97'
$ws.Range('E42').Value = '  +2.37%  '
$ws.Range('D43').Value = '''149.13'
$ws.Range('E43').Value = '  +4.07%  '
$ws.Range('E44').Value = '  +0.35%  '
$ws.Range('D45').Value = '''20.42'
$ws.Range('E45').Value = '  +1.14%  '
$ws.Range('E46').Value = '  +2.45%  '
$ws.Range('D47').Value = '''0.0959'
$ws.Range('E47').Value = '  -0.49%  '
$ws.Range('D48').Value = '''0.0513'
$ws.Range('E48').Value = '  -1.01%  '
$ws.Range('E49').Value = '  +3.25%  '
$ws.Range('D50').Value = '0.0₆0229'
$ws.Range('E50').Value = '  +4.81%  '
$ws.Range('D51').Value = '''17.83'
$ws.Range('E51').Value = '  +1.01%  '
